$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the sheet from 31 to 38 data rows by inserting 7 blank rows at the end
# of the existing table (before the old row 32 boundary, i.e. rows 32-38).
$ws.Range("A32:A38").EntireRow.Insert()

# Now (re)write every cell so the final table matches the target layout exactly,
# including rows that moved position and the newly-added rows.

# Row 1
$ws.Range("A1").Value = 'Location'
$ws.Range("B1").Value = 'Site'
$ws.Range("C1").Value = 'Exposure period'
$ws.Range("D1").Value = 'Notes'

# Row 2
$ws.Range("A2").Value = 'Brandon Park'
$ws.Range("B2").Value = 'Kmart, Brandon Park Shopping Centre  Cnr Springvale Rd and Ferntree Gully Rd  Brandon Park, VIC 3170'
$ws.Range("C2").Value = '4:35pm - 5:10pm  31/1/2021'
$ws.Range("D2").Value = 'Case attended venue'

# Row 3
$ws.Range("A3").Value = 'Brighton'
$ws.Range("B3").Value = 'North Point Cafe  2B North Rd  Brighton, VIC 3186'
$ws.Range("C3").Value = '8:10am - 9:30am  31/1/2021'
$ws.Range("D3").Value = 'Case dined outside and used bathroom'

# Row 4
$ws.Range("A4").Value = 'Broadmeadows'
$ws.Range("B4").Value = 'Craigieburn Line train'
$ws.Range("C4").Value = '1.25pm - 1.59pm  9/02/2012'
$ws.Range("D4").Value = 'Case caught train from Broadmeadows Railway Station to Glenroy Railway Station'

# Row 5
$ws.Range("A5").Value = 'Clayton South'
$ws.Range("B5").Value = 'Nakama Workshop  85 Main Rd  Clayton South, VIC 3169'
$ws.Range("C5").Value = '11:15am - 12:00pm  1/2/2021'
$ws.Range("D5").Value = 'Case attended venue'

# Row 6
$ws.Range("A6").Value = 'Coburg'
$ws.Range("B6").Value = 'Function venue  426 Sydney Rd  Coburg VIC 3058'
$ws.Range("C6").Value = '7:14pm  11:30pm  6/02/2021'
$ws.Range("D6").Value = 'Case attended venue'

# Row 7
$ws.Range("A7").Value = 'Glen Waverley'
$ws.Range("B7").Value = 'Commonwealth Bank, 28-32 Kingsway, Glen Waverley'
$ws.Range("C7").Value = '1:30pm-2:30pm 9/2/2021'
$ws.Range("D7").Value = 'Case attended venue'

# Row 8
$ws.Range("A8").Value = 'Glen Waverley'
$ws.Range("B8").Value = 'HSBC Bank, 38 Kingsway, Glen Waverley'
$ws.Range("C8").Value = '2:15pm-3:30pm 9/2/2021'
$ws.Range("D8").Value = 'Case attended venue'

# Row 9
$ws.Range("A9").Value = 'Glenroy'
$ws.Range("B9").Value = '513 Eltham to Glenroy bus route  Glenroy Railway Station towards Eltham'
$ws.Range("C9").Value = '1.35pm  2.17pm  9/02/2021'
$ws.Range("D9").Value = 'Case caught bus from Glenroy Railway Station towards Eltham'

# Row 10
$ws.Range("A10").Value = 'Heatherton'
$ws.Range("B10").Value = 'Melbourne Golf Academy  385 Centre Dandenong Rd  Heatherton, VIC 3202'
$ws.Range("C10").Value = '5:19pm - 6:35pm  1/2/2021'
$ws.Range("D10").Value = 'Case attended venue'

# Row 11
$ws.Range("A11").Value = 'Hoppers Crossing'
$ws.Range("B11").Value = 'Caltex Woolworths  50 Old Geelong Rd  Hoppers Crossing, VIC 3029'
$ws.Range("C11").Value = '6.40am - 7.15am  8/02/21'
$ws.Range("D11").Value = 'Case attended venue'

# Row 12
$ws.Range("A12").Value = 'Hoppers Crossing'
$ws.Range("B12").Value = 'Coates Hire Werribee  148A Geelong Rd  Hoppers Crossing, VIC 3029'
$ws.Range("C12").Value = '6.45am - 7.30am  8/02/21'
$ws.Range("D12").Value = 'Case attended venue'

# Row 13
$ws.Range("A13").Value = 'Keysborough'
$ws.Range("B13").Value = 'Aces Sporting Club (Driving Range)  Cnr Springvale Rd and Hutton Rd  Keysborough, VIC 3173'
$ws.Range("C13").Value = '10:00pm - 11:15pm  30/1/2021'
$ws.Range("D13").Value = 'Case attended venue'

# Row 14
$ws.Range("A14").Value = 'Keysborough'
$ws.Range("B14").Value = 'Kmart, Parkmore Keysborough Shopping Centre  C/317 Cheltenham Rd  Keysborough, VIC 3173'
$ws.Range("C14").Value = '4:00pm - 5:00pm  31/1/2021'
$ws.Range("D14").Value = 'Case attended venue'

# Row 15
$ws.Range("A15").Value = 'Maidstone'
$ws.Range("B15").Value = 'Marciano''s Cakes  126 Mitchell St  Maidstone VIC 3012'
$ws.Range("C15").Value = '9:45am - 10:25am  5/2/2021'
$ws.Range("D15").Value = 'Case attended venue'

# Row 16
$ws.Range("A16").Value = 'Melbourne'
$ws.Range("B16").Value = '901 Frankston to Melbourne Airport bus route  Melbourne Airport to Broadmeadows Railway Station'
$ws.Range("C16").Value = '1:02pm  1:49pm  9/2/2021'
$ws.Range("D16").Value = 'Case caught but from Melbourne Airport to Broadmeadows Railway Station'

# Row 17
$ws.Range("A17").Value = 'Melbourne'
$ws.Range("B17").Value = 'Brunetti: Terminal 4, Melbourne Airport'
$ws.Range("C17").Value = '4:45am - 1:15pm  9/2/2021'
$ws.Range("D17").Value = 'Case attended venue'

# Row 18
$ws.Range("A18").Value = 'Melbourne'
$ws.Range("B18").Value = 'Exford Hotel  199 Russell St  Melbourne, VIC 3000'
$ws.Range("C18").Value = '11:00pm - 11:35pm  29/1/2021'
$ws.Range("D18").Value = 'Case attended bottle shop'

# Row 19
$ws.Range("A19").Value = 'Melbourne'
$ws.Range("B19").Value = 'Terminal 4, Melbourne Airport'
$ws.Range("C19").Value = '4:45am - 2:00pm  9/2/2021'
$ws.Range("D19").Value = 'Case attended venue'

# Row 20
$ws.Range("A20").Value = 'Moorabbin Airport'
$ws.Range("B20").Value = 'Lululemon, DFO Moorabbin  Shop G-039/250 Centre Dandenong Rd  Moorabbin VIC 3194'
$ws.Range("C20").Value = '5:00pm - 5:45pm  1/2/2021'
$ws.Range("D20").Value = 'Case attended venue'

# Row 21
$ws.Range("A21").Value = 'Noble Park'
$ws.Range("B21").Value = 'Club Noble  46/56 Moodemere St  Noble Park VIC 3174'
$ws.Range("C21").Value = '2:36pm -3:30pm  30/01/2021'
$ws.Range("D21").Value = 'Case attended venue'

# Row 22
$ws.Range("A22").Value = 'South Melbourne'
$ws.Range("B22").Value = 'Stowe Australia  67  69 Buckhurst St  South Melbourne VIC 3205'
$ws.Range("C22").Value = '10.30am - 10.45am 8/02/2021'
$ws.Range("D22").Value = ""

# Row 23
$ws.Range("A23").Value = 'Springvale'
$ws.Range("B23").Value = 'Bunnings Springvale  849 Princes Hwy  Springvale, VIC 3171'
$ws.Range("C23").Value = '11:30am - 12:15pm  1/2/2021'
$ws.Range("D23").Value = 'Case attended venue'

# Row 24
$ws.Range("A24").Value = 'Springvale'
$ws.Range("B24").Value = 'Coles Springvale  825 Dandenong Rd  Springvale, VIC 3171'
$ws.Range("C24").Value = '5:00pm - 6:00pm  31/1/2021'
$ws.Range("D24").Value = 'Case attended venue'

# Row 25
$ws.Range("A25").Value = 'Springvale'
$ws.Range("B25").Value = 'Sharetea Springvale  27C Buckingham Ave  Springvale, VIC 3171'
$ws.Range("C25").Value = '6:50pm - 7:30pm  1/2/2021'
$ws.Range("D25").Value = 'Case attended venue'

# Row 26
$ws.Range("A26").Value = 'Springvale'
$ws.Range("B26").Value = 'Woolworths Springvale  302 Springvale Rd  Springvale, VIC 3171'
$ws.Range("C26").Value = '6:30pm - 7:30pm  1/2/2021'
$ws.Range("D26").Value = 'Case attended venue'

# Row 27
$ws.Range("A27").Value = 'Sunbury'
$ws.Range("B27").Value = 'Aldente Deli - Sunbury Square Shopping Centre  2-28 Evans Street  Sunbury VIC 3429'
$ws.Range("C27").Value = '3:45pm - 4:23pm  5/2/2021'
$ws.Range("D27").Value = 'Case attended venue'

# Row 28
$ws.Range("A28").Value = 'Sunbury'
$ws.Range("B28").Value = 'Asian Star - Sunbury Square Shopping Centre  2-28 Evans Street  Sunbury VIC 3429'
$ws.Range("C28").Value = '3:57pm - 4:30pm  5/2/2021'
$ws.Range("D28").Value = 'Case attended venue'

# Row 29
$ws.Range("A29").Value = 'Sunbury'
$ws.Range("B29").Value = 'Bakers Delight - Sunbury Square Shopping Centre  2-28 Evans Street  Sunbury VIC 3429'
$ws.Range("C29").Value = '3:40pm - 4:15pm  5/2/2021'
$ws.Range("D29").Value = 'Case attended venue'

# Row 30
$ws.Range("A30").Value = 'Sunbury'
$ws.Range("B30").Value = 'Cellarbrations  34 Batman Avenue  Sunbury VIC 3429'
$ws.Range("C30").Value = '5:44pm - 6:19pm  7/2/2021'
$ws.Range("D30").Value = 'Case attended venue'

# Row 31
$ws.Range("A31").Value = 'Sunbury'
$ws.Range("B31").Value = 'Cellarbrations  34 Batman Avenue  Sunbury VIC 3429'
$ws.Range("C31").Value = '6:17pm - 7:02pm  6/2/2021'
$ws.Range("D31").Value = 'Case attended venue'

# Row 32
$ws.Range("A32").Value = 'Sunbury'
$ws.Range("B32").Value = 'PJ''s Pet Warehouse  Shop 2, 104 Horne Street  Sunbury VIC 3429'
$ws.Range("C32").Value = '3:37pm - 4:10pm  5/2/2021'
$ws.Range("D32").Value = 'Case attended venue'

# Row 33
$ws.Range("A33").Value = 'Sunbury'
$ws.Range("B33").Value = 'Sunny Life Massage - Sunbury Square Shopping Centre  2-28 Evans Street  Sunbury VIC 3429'
$ws.Range("C33").Value = '4:30pm - 6:30pm  6/2/2021'
$ws.Range("D33").Value = 'Case attended venue'

# Row 34
$ws.Range("A34").Value = 'Sunbury'
$ws.Range("B34").Value = 'Sushi Sushi - Sunbury Square Shopping Centre  2-28 Evans Street  Sunbury VIC 3429'
$ws.Range("C34").Value = '3:53pm - 4:28pm  5/2/2021'
$ws.Range("D34").Value = 'Case attended venue'

# Row 35
$ws.Range("A35").Value = 'Sunshine'
$ws.Range("B35").Value = 'Dan Murphy''s  47 McIntyre Rd  Sunshine VIC 3020'
$ws.Range("C35").Value = '5:50pm - 6:30pm  5/2/2021'
$ws.Range("D35").Value = 'Case attended venue'

# Row 36
$ws.Range("A36").Value = 'Sunshine'
$ws.Range("B36").Value = 'Dan Murphy''s  47 McIntyre Rd  Sunshine VIC 3020'
$ws.Range("C36").Value = '6:50pm - 7:30pm  6/2/2021'
$ws.Range("D36").Value = 'Case attended venue'

# Row 37
$ws.Range("A37").Value = 'Taylors Lakes'
$ws.Range("B37").Value = 'Off Ya Tree Watergardens  399 Melton Highway  Taylors Lakes VIC 3038'
$ws.Range("C37").Value = '1:00pm - 1:52pm  6/2/2021'
$ws.Range("D37").Value = 'Case attended venue'

# Row 38
$ws.Range("A38").Value = 'West Melbourne'
$ws.Range("B38").Value = 'Kebab Kingz  438 Spencer St  West Melbourne, VIC 3003'
$ws.Range("C38").Value = '11:24pm -12:15am  29/1/2021'
$ws.Range("D38").Value = 'Case dined outside'
